# Update the public epexspot_prices.xlsx workbook:
#  1. "Prix Spot" sheet: insert a new date column before DL ("01-oct." ... shift
#     right to make room) for the missing "09-nov" day. New column gets "09-nov"
#     in the header row and "-" placeholders for every hourly data row.
#  2. "Gaz" and "CO2" sheets: append a new daily row (2025-11-07) after the
#     existing last row (144).

$wb = $excel.ActiveWorkbook

# --- 1. Prix Spot: insert missing "09-nov" day column at DL -----------------
$prix = $wb.Worksheets.Item("Prix Spot")

$prix.Range("DL1").EntireColumn.Insert()

$prix.Range("DL1").Value = "09-nov"
$prix.Range("DL2:DL25").Value = "-"

# --- 2. Gaz: append row 145 (2025-11-07) -------------------------------------
$gaz = $wb.Worksheets.Item("Gaz")

$gaz.Range("A145").NumberFormat = "@"
$gaz.Range("A145").Value = "2025-11-07"
$gaz.Range("A145").ClearFormats()
$gaz.Range("B145").Value = 29.74

# --- 3. CO2: append row 145 (2025-11-07) -------------------------------------
$co2 = $wb.Worksheets.Item("CO2")

$co2.Range("A145").NumberFormat = "@"
$co2.Range("A145").Value = "2025-11-07"
$co2.Range("A145").ClearFormats()
$co2.Range("B145").Value = 79.36
